# Java-02.pptx: drop the "Kurz Java 1" intro/WIFI slide, the "Termíny lekcí"
# schedule slide and the trailing "Objektové programování" slide, keeping only
# the title slide and the "Organizační pokyny" slide. Also strip the
# sponsor/partner logo pictures (ESF, OPZ, JetBrains) that were added to every
# slide layout, since they are no longer needed once the branded slides are
# gone.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# Remove the sponsor-logo picture shapes from every slide layout (the last
# shape(s) on each layout are the injected logo pictures). Delete from the
# highest shape index down so indices stay valid within each layout.
$master.CustomLayouts.Item(1).Shapes.Item(3).Delete()          # logoesfcrnatmavem.png
$master.CustomLayouts.Item(2).Shapes.Item(3).Delete()          # logoesfcrnatmavem.png
$master.CustomLayouts.Item(3).Shapes.Item(4).Delete()          # JetBrains picture
$master.CustomLayouts.Item(3).Shapes.Item(3).Delete()          # logoesfcrnatmavem.png
$master.CustomLayouts.Item(4).Shapes.Item(5).Delete()          # JetBrains picture
$master.CustomLayouts.Item(4).Shapes.Item(4).Delete()          # Logo-OPZ
$master.CustomLayouts.Item(5).Shapes.Item(4).Delete()          # Logo-OPZ
$master.CustomLayouts.Item(6).Shapes.Item(3).Delete()          # Logo-OPZ
$master.CustomLayouts.Item(7).Shapes.Item(4).Delete()          # Logo-OPZ

# Remove slides 5 ("Objektové programování"), 3 ("Termíny lekcí") and 2
# ("Kurz Java 1") — delete from the end so earlier indices remain valid.
$p.Slides.Item(5).Delete()
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()
